# Apply updated crypto price/volume data to Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number (e.g. "211.40")
# must keep their original Text storage (the sheet uses "." as a thousands
# separator in other rows, so these columns are authored as text, not numbers).
# Force Text format first so Excel does not reinterpret/round them as floats.
$textCells = @('D5', 'D8', 'D10', 'D11', 'D16', 'D19', 'D21', 'D22', 'D29', 'D38', 'D44', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.686.81'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '1.597.35'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '211.40'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.0618'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = '19.74'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').Value = '0.0839'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '1.821.64'
$ws.Range('E12').Value = '  -1.29%  '
$ws.Range('D13').Value = '1.626.44'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '65.11'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '26.696.01'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '210.34'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '4.27'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -3.84%  '
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = '15.31'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('E33').Value = '  -5.27%  '
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').Value = '1.296.75'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  -4.56%  '
$ws.Range('D38').Value = '0.0172'
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('E39').Value = '  +1.28%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('D44').Value = '63.61'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').Value = '1.733.62'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').Value = '90.15'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').Value = '0.872'
$ws.Range('E47').Value = '  +7.48%  '
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').Value = '0.0993'
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('D51').Value = '7.47'
$ws.Range('E51').Value = '  -0.59%  '
